$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain numeric-looking strings (e.g. "577.29") that must
# remain text cells (matching the source inlineStr cells) rather than being
# auto-converted to floating point numbers by Excel, which would both change the
# cell type and introduce binary floating point rounding noise. Forcing the
# NumberFormat to Text ("@") before the assignment keeps the literal string, and
# resetting the style to "Normal" afterwards avoids leaving a stray explicit
# number-format style on the cell (it reverts to the default style, as before).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.681.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.443.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.443.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.890.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.587.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("B18").Value = "BabyDogeCoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₅0119"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +328.39%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.461.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "328.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.07%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "639.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.18%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.00%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0980"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.72%  "
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("E35").Value = "  +4.05%  "
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +27.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.603"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.64%  "
